$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.795.08"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.799.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.81%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.799.64"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.73%  "

# Row 8
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$ws.Range("E9").Value = "  -0.55%  "

# Row 10
$ws.Range("E10").Value = "  -0.96%  "

# Row 11
$ws.Range("E11").Value = "  -0.91%  "

# Row 12
$ws.Range("E12").Value = "  -1.34%  "

# Row 13
$ws.Range("E13").Value = "  -3.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.433.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.802.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.56"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.59%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.741.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.41%  "

# Row 19
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("E20").Value = "  -0.24%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "461.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.38%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -9.18%  "

# Row 23
$ws.Range("E23").Value = "  -0.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000152"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.49"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.18%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.08"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.58%  "

# Row 27
$ws.Range("E27").Value = "  -4.42%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.947.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.80%  "

# Row 31
$ws.Range("E31").Value = "  -0.18%  "

# Row 32
$ws.Range("E32").Value = "  +2.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.68"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.05"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.91%  "

# Row 37
$ws.Range("E37").Value = "  -1.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.23%  "

# Row 39
$ws.Range("E39").Value = "  -0.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.58%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.94%  "

# Row 42
$ws.Range("E42").Value = "  +0.12%  "

# Row 43
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.22"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.62%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.77"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.61%  "

# Row 46
$ws.Range("E46").Value = "  -2.58%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.32"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.82%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.04%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "390.28"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.76%  "

# Row 51
$ws.Range("E51").Value = "  -6.16%  "
